$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.108930349349976
$ws.Range("B1").Value = 2.241667985916138
$ws.Range("C1").Value = 10.07673645019531
$ws.Range("D1").Value = 1.399994015693665
$ws.Range("E1").Value = 1.28227424621582
